$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 476, shifting existing rows 476:496 down to 477:497
$ws.Rows.Item(476).Insert()

# Populate the new row 476 with the new data record
$ws.Cells.Item(476, 1).Value = 4
$ws.Cells.Item(476, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(476, 3).Value = "Los Lagos"
$ws.Cells.Item(476, 4).Value = 45147
$ws.Cells.Item(476, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(476, 5).Value = 10
$ws.Cells.Item(476, 6).Value = 100112037
$ws.Cells.Item(476, 7).Value = "Cebollín"
$ws.Cells.Item(476, 8).Value = "Sin especificar"
$ws.Cells.Item(476, 9).Value = "Primera"
$ws.Cells.Item(476, 10).Value = 35
$ws.Cells.Item(476, 11).Value = 6000
$ws.Cells.Item(476, 12).Value = 6000
$ws.Cells.Item(476, 13).Value = 6000
$ws.Cells.Item(476, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(476, 15).Value = "Región Metropolitana"
$ws.Cells.Item(476, 16).Value = 167
$ws.Cells.Item(476, 17).Value = 36
$ws.Cells.Item(476, 18).Value = "Hortaliza"
